$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression (label unchanged)
$ws.Range("B2").Value = 0.1087152181300004
$ws.Range("C2").Value = 0.1087152181300004
$ws.Range("D2").Value = 0.1087152181300004

# Row 3: RandomForestRegressor (label unchanged)
$ws.Range("B3").Value = 0.02249433173081758
$ws.Range("C3").Value = 0.02131799511783486
$ws.Range("D3").Value = 0.0222990410571723

# Row 4: label changed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.0228165364663289
$ws.Range("C4").Value = 0.02178502945052628
$ws.Range("D4").Value = 0.02284553545631869

# Row 5: label changed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01674105237858478
$ws.Range("C5").Value = 0.0181454539629752
$ws.Range("D5").Value = 0.01840740124775413
